$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "合富中国"
$ws.Range("B2").Value = "平潭发展"
$ws.Range("C2").Value = "大鹏工业"
$ws.Range("A3").Value = "视觉中国"
$ws.Range("B3").Value = "合富中国"
$ws.Range("C3").Value = "平潭发展"
$ws.Range("A4").Value = "易点天下"
$ws.Range("B4").Value = "N大鹏"
$ws.Range("C4").Value = "航天发展"
$ws.Range("A5").Value = "航天发展"
$ws.Range("B5").Value = "航天发展"
$ws.Range("C5").Value = "榕基软件"
$ws.Range("A6").Value = "榕基软件"
$ws.Range("B6").Value = "视觉中国"
$ws.Range("C6").Value = "合富中国"
$ws.Range("A7").Value = "平潭发展"
$ws.Range("B7").Value = "榕基软件"
$ws.Range("C7").Value = "中水渔业"
$ws.Range("A8").Value = "国风新材"
$ws.Range("B8").Value = "易点天下"
$ws.Range("C8").Value = "视觉中国"
$ws.Range("A9").Value = "实达集团"
$ws.Range("B9").Value = "海南海药"
$ws.Range("C9").Value = "实达集团"
$ws.Range("A10").Value = "凯美特气"
$ws.Range("B10").Value = "海峡创新"
$ws.Range("C10").Value = "易点天下"
$ws.Range("A11").Value = "中水渔业"
$ws.Range("B11").Value = "九牧王"
$ws.Range("C11").Value = "国风新材"
$ws.Range("A12").Value = "海峡创新"
$ws.Range("B12").Value = "中水渔业"
$ws.Range("C12").Value = "九牧王"
$ws.Range("A13").Value = "九牧王"
$ws.Range("B13").Value = "浪潮软件"
$ws.Range("C13").Value = "华胜天成"
$ws.Range("A14").Value = "海南海药"
$ws.Range("B14").Value = "凯美特气"
$ws.Range("C14").Value = "海峡创新"
$ws.Range("A15").Value = "浪潮软件"
$ws.Range("B15").Value = "华映科技"
$ws.Range("C15").Value = "浪潮软件"
$ws.Range("A16").Value = "三花智控"
$ws.Range("B16").Value = "江龙船艇"
$ws.Range("C16").Value = "久其软件"
$ws.Range("A17").Value = "江龙船艇"
$ws.Range("B17").Value = "国风新材"
$ws.Range("C17").Value = "凯美特气"
$ws.Range("A18").Value = "华胜天成"
$ws.Range("B18").Value = "三花智控"
$ws.Range("C18").Value = "海南海药"
$ws.Range("A19").Value = ""
$ws.Range("B19").Value = "华胜天成"
$ws.Range("C19").Value = "华夏幸福"
$ws.Range("A20").Value = ""
$ws.Range("B20").Value = "实达集团"
$ws.Range("C20").Value = "江龙船艇"
$ws.Range("A21").Value = ""
$ws.Range("B21").Value = "大众公用"
$ws.Range("C21").Value = "海马汽车"
